$d = $word.ActiveDocument

# wdReplaceAll = 2
$replaceAll = 2

# 1) Title / "Θέμα" sentence:
#    «${school}» στη χώρα «${country}»,  στο πλαίσιο ...
# -> «${school}» ${country},  στο πλαίσιο ...
$d.Content.Find.Execute(' στη χώρα «${country}»', $true, $false, $false, $false, $false, $true, 1, $false, ' ${country}', $replaceAll)

# 2) Body sentence:
#    ... εταιρικό σχολείο στη «${country}» από ...
# -> ... εταιρικό σχολείο ${country} από ...
$d.Content.Find.Execute('εταιρικό σχολείο στη «${country}»', $true, $false, $false, $false, $false, $true, 1, $false, 'εταιρικό σχολείο ${country}', $replaceAll)
